# Updates the "Price" (column D) and a couple of "Volume(1h)" (column E)
# values on Sheet1 to match the refreshed symbol-list snapshot.
#
# The source cells are stored as plain TEXT (not numbers), even though most
# of them look numeric (e.g. "248.77"). A bare
#   $ws.Range("D2").Value = "248.72"
# would let Excel's smart-parsing turn that into a genuine number, which
# would change the cell's stored type. To keep these as text, values are
# assigned with a leading apostrophe (forces literal-text interpretation,
# same as typing '248.72 into a cell), and then the cell's Style is reset
# to "Normal" so no stray "quote prefix" / text-number-format styling is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $ws.Range($cellRef).Value = "'" + $newValue
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2"  "248.72"
Set-TextValue "D3"  "22.53"
Set-TextValue "D4"  "5.430"
Set-TextValue "D5"  "0.05695"
Set-TextValue "D6"  "3.393"
Set-TextValue "D7"  "6.319"
Set-TextValue "D8"  "0.8125"
Set-TextValue "D9"  "0.9319"
Set-TextValue "D10" "0.1415"
Set-TextValue "D11" "0.07420"
Set-TextValue "D12" "0.03042"
Set-TextValue "D13" "0.03017"
Set-TextValue "D14" "0.09375"
Set-TextValue "D15" "3.713"
Set-TextValue "D16" "0.001576"
Set-TextValue "D17" "0.04749"
Set-TextValue "E19" "18OneONE"
Set-TextValue "D20" "0.006440"
Set-TextValue "D21" "0.005000"
Set-TextValue "D24" "3.696"
Set-TextValue "D25" "2.163"
Set-TextValue "D40" "0.04004"
Set-TextValue "D41" "0.1067"
Set-TextValue "D43" "0.002997"
Set-TextValue "E43" "42KickTokenKICKWorstin24h"
Set-TextValue "D44" "0.007488"
Set-TextValue "D47" "0.4999"
